$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = 340926355
$ws.Range("C14").Value = 5985666
$ws.Range("C15").Value = 5530000000
$ws.Range("C16").Value = 51786608
$ws.Range("C18").Formula = "=SUM(C12:C17)"
$ws.Range("C19").Value = 392700000
$ws.Range("C21").Formula = "=SUM(C18:C20)"
$ws.Range("C22").Value = 3376095
$ws.Range("C26").Value = 1003368421

$excel.Calculate()
